$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# ---------------------------------------------------------------------------
# New query text for the "StatQuery" column (now shared by every tab row)
# ---------------------------------------------------------------------------
$statQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['Akita']`nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"

# ---------------------------------------------------------------------------
# New CasesTab query (adds a Cohort column + integer-safe Age expression)
# ---------------------------------------------------------------------------
$casesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Akita']`nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age`nRETURN  `n       coalesce(c.case_id, '') AS ``Case ID``,`n       coalesce(s.clinical_study_designation, '') AS ``Study Code``,`n       coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n       coalesce(demo.breed, '') AS Breed ,`n       coalesce(diag.disease_term, '') AS Diagnosis ,`n       coalesce(diag.stage_of_disease, '') AS ``Stage of Disease``,`n       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,`n       coalesce(demo.sex, '') AS Sex,`n       coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n       coalesce(demo.weight, '') AS ``Weight (kg)``,`n       coalesce(diag.best_response, '') AS ``Response to Treatment``,`n       coalesce(co.cohort_description, '') AS ``Cohort``"

# ---------------------------------------------------------------------------
# New FilesTab query (adds Sample ID, human-readable Size/Format columns)
# ---------------------------------------------------------------------------
$filesQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f)-[*]->(samp:sample)`n MATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Akita']`nOPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)`nWITH`n        f, parent, c, demo, diag, s, samp,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH    `n        f, parent, c, demo, diag, s, samp,`n        f.file_size /(1024^i) AS value, 10^precision AS factor,`n        units[i] as unit`nRETURN coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(f.file_type, '') AS ``File Type``,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        round(factor * value)/factor+' ' +unit AS Size,`n        coalesce(samp.sample_id, '') AS ``Sample ID``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(demo.breed,'') AS Breed ,`n        coalesce(diag.disease_term,'') AS Diagnosis"

# ---------------------------------------------------------------------------
# Brand new StudyFilesTab query
# ---------------------------------------------------------------------------
$studyFilesQuery = "  MATCH (f:file)-->(s:study)`n MATCH (f)-->(parent)`nMATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`nWHERE demo.breed  IN ['Akita']  `nWITH DISTINCT f, parent, s, c, demo, diag`nWITH`n        f, c, demo, diag, s,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH    `n        f, c, demo, diag, s,`n        f.file_size /(1024^i) AS value, 10^precision AS factor,`n        units[i] as unit`nRETURN `n  coalesce(f.file_name, '') AS ``File Name``,`n  coalesce(f.file_type, '') AS ``File Type``,`n  coalesce(""study"", '') AS ``Association``,`n  coalesce(f.file_description, '') AS ``Description``,`n  coalesce(f.file_format, '') AS  Format,`n  round(factor * value)/factor+' ' +unit AS Size,`n  coalesce(s.clinical_study_designation,'') AS ``Study Code``"

# ---------------------------------------------------------------------------
# Update the StatQuery column (C) for the three existing rows
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# Update the query column (B) for the CasesTab / FilesTab rows
$ws.Range("B2").Value = $casesQuery
$ws.Range("B4").Value = $filesQuery

# ---------------------------------------------------------------------------
# Drop the old trailing placeholder rows (6-13), keep row 5 for the new tab
# ---------------------------------------------------------------------------
$ws.Rows("6:13").Delete()

# ---------------------------------------------------------------------------
# Populate the new "StudyFilesTab" row (row 5)
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = $studyFilesQuery
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("C5").Value = $statQuery
$ws.Range("D5").Value = "TC01_Canine_Filter_Breed-Akita_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC01_Canine_Filter_Breed-Akita_WebData.xlsx"

$ws.Range("B5").WrapText = $true
$ws.Range("C5").WrapText = $true

# ---------------------------------------------------------------------------
# Row heights (Excel recomputed these once the wrapped text changed)
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 285
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 213

# ---------------------------------------------------------------------------
# View state: zoomed out, scrolled to the new row, new selection
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 55
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B5").Select()
